$wb = $excel.ActiveWorkbook

# Add a new worksheet named "About" after the last existing sheet so it
# lands at the end of the tab strip (and becomes the active sheet, matching
# the activeTab/tabSelected changes in the diff).
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$about = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$about.Name = "About"

# Header row
$about.Range("A1").Value = "country"
$about.Range("B1").Value = "about"
$about.Range("C1").Value = "prone_to"

# Data entered country-by-country in the same order the countries first
# appear elsewhere in the workbook (Ruritania, then Aurelia, then Xenon),
# even though the rows end up sorted differently on the sheet.
$about.Range("A3").Value = "Ruritania"
$about.Range("B3").Value = "Ruritania is a tropical country whose main export is agricultural products. Before the pandemic the government was steadily making efforts consolidating the fiscal accounts which translated into a rating upgrade in 2021. The pandemic deteriorated the accounts and since then the government has been working to shore them up."
$about.Range("C3").Value = "The country has been affected by droughts that have reduced its agricultural production. Heatwaves have decreased productivity and increased energy and water consumption. The country could benefit from better irrigation infrastructure and water management."

$about.Range("A2").Value = "Aurelia"
$about.Range("B2").Value = "Aurelia is on the coast and has been affected by large hurricanes recently. With heavy rains there have been some serious landslides that have destroyed roads and other infrastructure. Climate scientist project that these events will increase in magnitude and frequency. The government is looking for financing to protect its growth and fiscal accounts from future events."
$about.Range("C2").Value = "Rising sea levels have affected coastal towns and cities. Tourism is particularly vulnerable as beaches get narrower and infrastructure deteriorated. There are projects in place to improve protection for sea fronts and ports."

$about.Range("A4").Value = "Xenon"
$about.Range("B4").Value = "Xenon is a country that relies on labor-intensive manufacturing. The country has been affected by flooding due to unusually heavy rain. The government has built waterways and flood protection, but this has deteriorated its fiscal accounts in the process."
$about.Range("C4").Value = "The government is interested in finding solutions that could improve its resilience to weather events without deteriorating its fiscal accounts. Borrowing cost increased after its credit rating fell from BB- to B+ in 2024."

# Data rows use a new font (Aptos, 12pt, theme text color) distinct from
# the workbook default.
$about.Range("A2:C4").Font.Name = "Aptos"
$about.Range("A2:C4").Font.Size = 12

# Match the recorded selection on the new sheet.
$about.Range("D12").Select()
